# Update the weekly Fruta/Hortaliza price records.
# The underlying data rows (2, 3, 5) get re-ordered/re-valued as new
# weekly records arrive; row 4 is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- becomes the "Segunda" quality record (previously on row 5)
$ws.Range("D2").Value = 44881
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11250
$ws.Range("O2").Value = 11250
$ws.Range("P2").Value = 11250
$ws.Range("S2").Value = 11250

# Row 3 <- becomes the first "Primera" quality record (previously on row 2)
$ws.Range("D3").Value = 44874
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 7500
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 7750
$ws.Range("S3").Value = 7750

# Row 4 stays the same (already matches the target "Primera" record).

# Row 5 <- becomes the other "Primera" quality record (previously on rows 3/4)
$ws.Range("D5").Value = 44923
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 7500
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 7625
$ws.Range("S5").Value = 7625
